$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 400
$ws.Range("B2").Value = 396
$ws.Range("C2").Value = 0

$ws.Range("B5").Value = 0.99
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.01
